$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 45909
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B5").Value = "20,9478"
$ws.Range("C5").Value = "14,7352"
$ws.Range("D5").Value = "14,8086"
$ws.Range("E5").Value = "14,8086"
